# Update LR-pair (Inhbb-Acvr1b) NATMI stats with refreshed TPM-based numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.372558333333333
$ws.Cells.Item(2, 8).Value = 4.117675
$ws.Cells.Item(2, 9).Value = 0.3340102211301095
$ws.Cells.Item(2, 10).Value = 0.3340102211301095
$ws.Cells.Item(2, 13).Value = 3.483060666666667
$ws.Cells.Item(2, 14).Value = 10.449182
$ws.Cells.Item(2, 15).Value = 0.2527672867110271
$ws.Cells.Item(2, 16).Value = 0.2527672867110271
$ws.Cells.Item(2, 17).Value = 4.780703943538889
$ws.Cells.Item(2, 18).Value = 43.02633549185001
$ws.Cells.Item(2, 19).Value = 0.08442685732880795
$ws.Cells.Item(2, 20).Value = 0.08442685732880796

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.372558333333333
$ws.Cells.Item(3, 8).Value = 4.117675
$ws.Cells.Item(3, 9).Value = 0.3340102211301095
$ws.Cells.Item(3, 10).Value = 0.3340102211301095
$ws.Cells.Item(3, 15).Value = 0.3353267952677969
$ws.Cells.Item(3, 16).Value = 0.335326795267797
$ws.Cells.Item(3, 17).Value = 6.342189898741666
$ws.Cells.Item(3, 18).Value = 57.079709088675
$ws.Cells.Item(3, 19).Value = 0.1120025770382478
$ws.Cells.Item(3, 20).Value = 0.1120025770382478

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.372558333333333
$ws.Cells.Item(4, 8).Value = 4.117675
$ws.Cells.Item(4, 9).Value = 0.3340102211301095
$ws.Cells.Item(4, 10).Value = 0.3340102211301095
$ws.Cells.Item(4, 13).Value = 2.773309666666667
$ws.Cells.Item(4, 14).Value = 8.319929
$ws.Cells.Item(4, 15).Value = 0.2012603358768551
$ws.Cells.Item(4, 16).Value = 0.2012603358768551
$ws.Cells.Item(4, 17).Value = 3.806529293897222
$ws.Cells.Item(4, 18).Value = 34.258763645075
$ws.Cells.Item(4, 19).Value = 0.06722300929094849
$ws.Cells.Item(4, 20).Value = 0.0672230092909485

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.372558333333333
$ws.Cells.Item(5, 8).Value = 4.117675
$ws.Cells.Item(5, 9).Value = 0.3340102211301095
$ws.Cells.Item(5, 10).Value = 0.3340102211301095
$ws.Cells.Item(5, 13).Value = 2.902635666666666
$ws.Cells.Item(5, 14).Value = 8.707906999999999
$ws.Cells.Item(5, 15).Value = 0.2106455821443209
$ws.Cells.Item(5, 16).Value = 0.2106455821443209
$ws.Cells.Item(5, 17).Value = 3.984036772913889
$ws.Cells.Item(5, 18).Value = 35.85633095622499
$ws.Cells.Item(5, 19).Value = 0.07035777747210527
$ws.Cells.Item(5, 20).Value = 0.07035777747210527

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 9).Value = 0.01293592767872722
$ws.Cells.Item(6, 10).Value = 0.01293592767872721
$ws.Cells.Item(6, 13).Value = 3.483060666666667
$ws.Cells.Item(6, 14).Value = 10.449182
$ws.Cells.Item(6, 15).Value = 0.2527672867110271
$ws.Cells.Item(6, 16).Value = 0.2527672867110271
$ws.Cells.Item(6, 17).Value = 0.1851525389186667
$ws.Cells.Item(6, 18).Value = 1.666372850268
$ws.Cells.Item(6, 19).Value = 0.003269779340441953
$ws.Cells.Item(6, 20).Value = 0.003269779340441953

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 9).Value = 0.01293592767872722
$ws.Cells.Item(7, 10).Value = 0.01293592767872721
$ws.Cells.Item(7, 15).Value = 0.3353267952677969
$ws.Cells.Item(7, 16).Value = 0.335326795267797
$ws.Cells.Item(7, 19).Value = 0.004337763172323588
$ws.Cells.Item(7, 20).Value = 0.004337763172323589

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 9).Value = 0.01293592767872722
$ws.Cells.Item(8, 10).Value = 0.01293592767872721
$ws.Cells.Item(8, 13).Value = 2.773309666666667
$ws.Cells.Item(8, 14).Value = 8.319929
$ws.Cells.Item(8, 15).Value = 0.2012603358768551
$ws.Cells.Item(8, 16).Value = 0.2012603358768551
$ws.Cells.Item(8, 17).Value = 0.1474235952606667
$ws.Cells.Item(8, 18).Value = 1.326812357346
$ws.Cells.Item(8, 19).Value = 0.002603489149499346
$ws.Cells.Item(8, 20).Value = 0.002603489149499346

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 9).Value = 0.01293592767872722
$ws.Cells.Item(9, 10).Value = 0.01293592767872721
$ws.Cells.Item(9, 13).Value = 2.902635666666666
$ws.Cells.Item(9, 14).Value = 8.707906999999999
$ws.Cells.Item(9, 15).Value = 0.2106455821443209
$ws.Cells.Item(9, 16).Value = 0.2106455821443209
$ws.Cells.Item(9, 17).Value = 0.1542983067686667
$ws.Cells.Item(9, 18).Value = 1.388684760918
$ws.Cells.Item(9, 19).Value = 0.002724896016462328
$ws.Cells.Item(9, 20).Value = 0.002724896016462328

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 7).Value = 2.683614
$ws.Cells.Item(10, 8).Value = 8.050841999999999
$ws.Cells.Item(10, 9).Value = 0.6530538511911632
$ws.Cells.Item(10, 10).Value = 0.6530538511911632
$ws.Cells.Item(10, 13).Value = 3.483060666666667
$ws.Cells.Item(10, 14).Value = 10.449182
$ws.Cells.Item(10, 15).Value = 0.2527672867110271
$ws.Cells.Item(10, 16).Value = 0.2527672867110271
$ws.Cells.Item(10, 17).Value = 9.347190367915999
$ws.Cells.Item(10, 18).Value = 84.12471331124399
$ws.Cells.Item(10, 19).Value = 0.1650706500417772
$ws.Cells.Item(10, 20).Value = 0.1650706500417772

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 7).Value = 2.683614
$ws.Cells.Item(11, 8).Value = 8.050841999999999
$ws.Cells.Item(11, 9).Value = 0.6530538511911632
$ws.Cells.Item(11, 10).Value = 0.6530538511911632
$ws.Cells.Item(11, 15).Value = 0.3353267952677969
$ws.Cells.Item(11, 16).Value = 0.335326795267797
$ws.Cells.Item(11, 17).Value = 12.400193995098
$ws.Cells.Item(11, 18).Value = 111.601745955882
$ws.Cells.Item(11, 19).Value = 0.2189864550572255
$ws.Cells.Item(11, 20).Value = 0.2189864550572255

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 7).Value = 2.683614
$ws.Cells.Item(12, 8).Value = 8.050841999999999
$ws.Cells.Item(12, 9).Value = 0.6530538511911632
$ws.Cells.Item(12, 10).Value = 0.6530538511911632
$ws.Cells.Item(12, 13).Value = 2.773309666666667
$ws.Cells.Item(12, 14).Value = 8.319929
$ws.Cells.Item(12, 15).Value = 0.2012603358768551
$ws.Cells.Item(12, 16).Value = 0.2012603358768551
$ws.Cells.Item(12, 17).Value = 7.442492647801999
$ws.Cells.Item(12, 18).Value = 66.982433830218
$ws.Cells.Item(12, 19).Value = 0.1314338374364073
$ws.Cells.Item(12, 20).Value = 0.1314338374364073

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 7).Value = 2.683614
$ws.Cells.Item(13, 8).Value = 8.050841999999999
$ws.Cells.Item(13, 9).Value = 0.6530538511911632
$ws.Cells.Item(13, 10).Value = 0.6530538511911632
$ws.Cells.Item(13, 13).Value = 2.902635666666666
$ws.Cells.Item(13, 14).Value = 8.707906999999999
$ws.Cells.Item(13, 15).Value = 0.2106455821443209
$ws.Cells.Item(13, 16).Value = 0.2106455821443209
$ws.Cells.Item(13, 17).Value = 7.789553711965999
$ws.Cells.Item(13, 18).Value = 70.10598340769398
$ws.Cells.Item(13, 19).Value = 0.1375629086557533
$ws.Cells.Item(13, 20).Value = 0.1375629086557533
